$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix capitalization in the English footnote (C7's "*according" -> "*According")
$footnote = $ws.Range("C7").Value()
$ws.Range("C7").Value = $footnote.Replace("*according", "*According")

# 2. Add the new 2023 column (O) to the table, mirroring the existing 2022 column (N)
#    Header row (row 3)
$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("O3").Value = 2023

#    Data rows (4-6)
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value = 5571

$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").Value = 74710

$ws.Range("N6").Copy()
$ws.Range("O6").PasteSpecial(-4122)
$ws.Range("O6").Value = 375715

#    Also extend the thin divider/blank formatting in row 2 to column O
$ws.Range("N2").Copy()
$ws.Range("O2").PasteSpecial(-4122)

# 3. Resize the header row and the footnote row to better fit the wider table
$ws.Rows.Item(1).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 46.5

# 4. Shrink the footnote row (A7:C7) font from 9pt to 8pt so it fits the new layout
$ws.Range("A7").Font.Size = 8
$ws.Range("B7:C7").Font.Size = 8
